# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" rows (16-21) get re-sorted into ascending period order
# (1901 -> 1906) and the per-period "Valor Mora" amounts travel with their
# period label. The account-statement logo picture also shifts slightly to
# the left.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-sort the "Periodo Mora" (E16:E21) / "Valor Mora" (F16:F21) rows ---
# Before: 1906,1905,1904,1903,1902,1901 with values 17708,31249,31249,31249,31249,31249
# After:  1901,1902,1903,1904,1905,1906 with values 31249,31249,31249,31249,31249,17708
$periods = @("1901", "1902", "1903", "1904", "1905", "1906")
$values  = @(31249, 31249, 31249, 31249, 31249, 17708)

for ($i = 0; $i -lt 6; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
    $ws.Cells.Item($row, 6).Value = $values[$i]
}

# --- Nudge the logo image 13.5pt (171450 EMU) to the left ---
# (Using a precise absolute target rather than "$shp.Left - 13.5" since the
# Left getter rounds to 2dp for display and that rounding would otherwise
# bleed into the saved anchor offsets.)
if ($ws.Shapes.Count -ge 1) {
    $shp = $ws.Shapes.Item(1)
    $shp.Left = 61.912817
}
